$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (not number/date auto-detection) for columns I, Y, AA
# across the new rows, matching the source export which stores these as text.
$ws.Range("I48:I53").NumberFormat = "@"
$ws.Range("Y48:Y53").NumberFormat = "@"
$ws.Range("AA48:AA53").NumberFormat = "@"

# Row 48
$ws.Range("A48").Value = 111974126
$ws.Range("B48").Value = 88032
$ws.Range("C48").Value = 'Ovaliderad'
$ws.Range("D48").Value = 'VU'
$ws.Range("E48").Value = 6276
$ws.Range("F48").Value = 'Goliatmusseron'
$ws.Range("G48").Value = 'Tricholoma matsutake'
$ws.Range("H48").Value = '(S.Ito & S.Imai) Singer'
$ws.Range("I48").Value = '1'
$ws.Range("J48").Value = 'mycel'
$ws.Range("P48").Value = 'Aloppmoarna i S, Jmt'
$ws.Range("Q48").Value = 439289.9461055733
$ws.Range("R48").Value = 6952209.002200785
$ws.Range("S48").Value = 10
$ws.Range("T48").Value = 'Jämtland'
$ws.Range("U48").Value = 'Härjedalen'
$ws.Range("V48").Value = 'Jämtland'
$ws.Range("W48").Value = 'Vemdalen'
$ws.Range("Y48").Value = '2023-09-05'
$ws.Range("Z48").Value = '00:00'
$ws.Range("AA48").Value = '2023-09-05'
$ws.Range("AB48").Value = '00:00'
$ws.Range("AD48").Value = $false
$ws.Range("AE48").Value = $false
$ws.Range("AG48").Value = $false
$ws.Range("AI48").Value = 'äldre renbetad ristallskog med lavfläckar på torr moränmark'
# AT48: empty inlineStr cell in source — no COM equivalent (Value="" deletes the cell); left unset
$ws.Range("AW48").Value = 'Magnus Andersson'
$ws.Range("AX48").Value = 'Magnus Andersson'
$ws.Range("AY48").Value = 'SCA Skog Naturvärdesinventering'

# Row 49
$ws.Range("A49").Value = 111974124
$ws.Range("B49").Value = 90666
$ws.Range("C49").Value = 'Ovaliderad'
$ws.Range("D49").Value = 'LC'
$ws.Range("E49").Value = 4364
$ws.Range("F49").Value = 'Dropptaggsvamp'
$ws.Range("G49").Value = 'Hydnellum ferrugineum'
$ws.Range("H49").Value = '(Fr.:Fr.) P. Karst.'
$ws.Range("I49").Value = '1'
$ws.Range("J49").Value = 'mycel'
$ws.Range("P49").Value = 'Aloppmoarna i S, Jmt'
$ws.Range("Q49").Value = 439276.3867801811
$ws.Range("R49").Value = 6952196.853249942
$ws.Range("S49").Value = 10
$ws.Range("T49").Value = 'Jämtland'
$ws.Range("U49").Value = 'Härjedalen'
$ws.Range("V49").Value = 'Jämtland'
$ws.Range("W49").Value = 'Vemdalen'
$ws.Range("Y49").Value = '2023-09-05'
$ws.Range("Z49").Value = '00:00'
$ws.Range("AA49").Value = '2023-09-05'
$ws.Range("AB49").Value = '00:00'
$ws.Range("AD49").Value = $false
$ws.Range("AE49").Value = $false
$ws.Range("AG49").Value = $false
$ws.Range("AI49").Value = 'äldre renbetad ristallskog med lavfläckar på torr moränmark'
# AT49: empty inlineStr cell in source — no COM equivalent (Value="" deletes the cell); left unset
$ws.Range("AW49").Value = 'Magnus Andersson'
$ws.Range("AX49").Value = 'Magnus Andersson'
$ws.Range("AY49").Value = 'SCA Skog Naturvärdesinventering'

# Row 50
$ws.Range("A50").Value = 111974133
$ws.Range("B50").Value = 90682
$ws.Range("C50").Value = 'Ovaliderad'
$ws.Range("D50").Value = 'NT'
$ws.Range("E50").Value = 2059
$ws.Range("F50").Value = 'Skrovlig taggsvamp'
$ws.Range("G50").Value = 'Hydnellum scabrosum'
$ws.Range("H50").Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("I50").Value = '1'
$ws.Range("J50").Value = 'mycel'
$ws.Range("P50").Value = 'Aloppmoarna i S, Jmt'
$ws.Range("Q50").Value = 439389.9449806474
$ws.Range("R50").Value = 6952220.480550999
$ws.Range("S50").Value = 10
$ws.Range("T50").Value = 'Jämtland'
$ws.Range("U50").Value = 'Härjedalen'
$ws.Range("V50").Value = 'Jämtland'
$ws.Range("W50").Value = 'Vemdalen'
$ws.Range("Y50").Value = '2023-09-05'
$ws.Range("Z50").Value = '00:00'
$ws.Range("AA50").Value = '2023-09-05'
$ws.Range("AB50").Value = '00:00'
$ws.Range("AD50").Value = $false
$ws.Range("AE50").Value = $false
$ws.Range("AG50").Value = $false
$ws.Range("AI50").Value = 'äldre renbetad ristallskog med lavfläckar på torr moränmark'
# AT50: empty inlineStr cell in source — no COM equivalent (Value="" deletes the cell); left unset
$ws.Range("AW50").Value = 'Magnus Andersson'
$ws.Range("AX50").Value = 'Magnus Andersson'
$ws.Range("AY50").Value = 'SCA Skog Naturvärdesinventering'

# Row 51
$ws.Range("A51").Value = 111974029
$ws.Range("B51").Value = 88032
$ws.Range("C51").Value = 'Ovaliderad'
$ws.Range("D51").Value = 'VU'
$ws.Range("E51").Value = 6276
$ws.Range("F51").Value = 'Goliatmusseron'
$ws.Range("G51").Value = 'Tricholoma matsutake'
$ws.Range("H51").Value = '(S.Ito & S.Imai) Singer'
$ws.Range("I51").Value = '1'
$ws.Range("J51").Value = 'mycel'
$ws.Range("P51").Value = 'Aloppmoarna, Jmt'
$ws.Range("Q51").Value = 439334.7866423383
$ws.Range("R51").Value = 6952296.802153576
$ws.Range("S51").Value = 10
$ws.Range("T51").Value = 'Jämtland'
$ws.Range("U51").Value = 'Härjedalen'
$ws.Range("V51").Value = 'Jämtland'
$ws.Range("W51").Value = 'Vemdalen'
$ws.Range("Y51").Value = '2023-09-05'
$ws.Range("Z51").Value = '00:00'
$ws.Range("AA51").Value = '2023-09-05'
$ws.Range("AB51").Value = '00:00'
$ws.Range("AD51").Value = $false
$ws.Range("AE51").Value = $false
$ws.Range("AG51").Value = $false
$ws.Range("AI51").Value = 'äldre renbetad ristallskog med lavfläckar på torr moränmark'
# AT51: empty inlineStr cell in source — no COM equivalent (Value="" deletes the cell); left unset
$ws.Range("AW51").Value = 'Magnus Andersson'
$ws.Range("AX51").Value = 'Magnus Andersson'
$ws.Range("AY51").Value = 'SCA Skog Naturvärdesinventering'

# Row 52
$ws.Range("A52").Value = 111974125
$ws.Range("B52").Value = 90660
$ws.Range("C52").Value = 'Ovaliderad'
$ws.Range("D52").Value = 'NT'
$ws.Range("E52").Value = 4362
$ws.Range("F52").Value = 'Blå taggsvamp'
$ws.Range("G52").Value = 'Hydnellum caeruleum'
$ws.Range("H52").Value = '(Hornem.) P.Karst.'
$ws.Range("I52").Value = '1'
$ws.Range("J52").Value = 'mycel'
$ws.Range("P52").Value = 'Aloppmoarna i S, Jmt'
$ws.Range("Q52").Value = 439278.8711310769
$ws.Range("R52").Value = 6952206.909989387
$ws.Range("S52").Value = 10
$ws.Range("T52").Value = 'Jämtland'
$ws.Range("U52").Value = 'Härjedalen'
$ws.Range("V52").Value = 'Jämtland'
$ws.Range("W52").Value = 'Vemdalen'
$ws.Range("Y52").Value = '2023-09-05'
$ws.Range("Z52").Value = '00:00'
$ws.Range("AA52").Value = '2023-09-05'
$ws.Range("AB52").Value = '00:00'
$ws.Range("AD52").Value = $false
$ws.Range("AE52").Value = $false
$ws.Range("AG52").Value = $false
$ws.Range("AI52").Value = 'äldre renbetad ristallskog med lavfläckar på torr moränmark'
# AT52: empty inlineStr cell in source — no COM equivalent (Value="" deletes the cell); left unset
$ws.Range("AW52").Value = 'Magnus Andersson'
$ws.Range("AX52").Value = 'Magnus Andersson'
$ws.Range("AY52").Value = 'SCA Skog Naturvärdesinventering'

# Row 53
$ws.Range("A53").Value = 111974134
$ws.Range("B53").Value = 90658
$ws.Range("C53").Value = 'Ovaliderad'
$ws.Range("D53").Value = 'NT'
$ws.Range("E53").Value = 4361
$ws.Range("F53").Value = 'Orange taggsvamp'
$ws.Range("G53").Value = 'Hydnellum aurantiacum'
$ws.Range("H53").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("I53").Value = '1'
$ws.Range("J53").Value = 'mycel'
$ws.Range("P53").Value = 'Aloppmoarna i S, Jmt'
$ws.Range("Q53").Value = 439399.8222122483
$ws.Range("R53").Value = 6952207.441512506
$ws.Range("S53").Value = 10
$ws.Range("T53").Value = 'Jämtland'
$ws.Range("U53").Value = 'Härjedalen'
$ws.Range("V53").Value = 'Jämtland'
$ws.Range("W53").Value = 'Vemdalen'
$ws.Range("Y53").Value = '2023-09-05'
$ws.Range("Z53").Value = '00:00'
$ws.Range("AA53").Value = '2023-09-05'
$ws.Range("AB53").Value = '00:00'
$ws.Range("AD53").Value = $false
$ws.Range("AE53").Value = $false
$ws.Range("AG53").Value = $false
$ws.Range("AI53").Value = 'äldre renbetad ristallskog med lavfläckar på torr moränmark'
# AT53: empty inlineStr cell in source — no COM equivalent (Value="" deletes the cell); left unset
$ws.Range("AW53").Value = 'Magnus Andersson'
$ws.Range("AX53").Value = 'Magnus Andersson'
$ws.Range("AY53").Value = 'SCA Skog Naturvärdesinventering'
